$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a new row after 47 by copying row 47 entirely (shifts everything below)
$ws.Rows.Item(47).Copy()
$ws.Rows.Item(48).Insert()

# Step 2: Make row 47 look like a normal (non-last) row by copying formatting from row 46 (only used range)
$ws.Range("B46:J46").Copy()
$ws.Range("B47:J47").PasteSpecial(-4122)  # xlPasteFormats

Write-Output "done"
